$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 692.2308
$ws.Range("I12").Value = 499.9091
$ws.Range("K12").Value = 499.9091
$ws.Range("M12").Value = -329.9091
$ws.Range("H17").Value = 607.2778
$ws.Range("I17").Value = 663
$ws.Range("J17").Value = 604.8551
$ws.Range("K17").Value = 1989
$ws.Range("L17").Value = 1814.5653
$ws.Range("M17").Value = -1821
$ws.Range("N17").Value = -2150.5653
$ws.Range("H28").Value = 9387.143
$ws.Range("J28").Value = 2798.5
$ws.Range("L28").Value = 2798.5
$ws.Range("N28").Value = -3768.5
$ws.Range("H48").Value = 5122.324
$ws.Range("I48").Value = 2000
$ws.Range("J48").Value = 5209.0557
$ws.Range("K48").Value = 6000
$ws.Range("L48").Value = 15627.1671
$ws.Range("M48").Value = -5708
$ws.Range("N48").Value = -16211.1671
$ws.Range("H56").Value = 5122.324
$ws.Range("I56").Value = 2000
$ws.Range("J56").Value = 5209.0557
$ws.Range("K56").Value = 6000
$ws.Range("L56").Value = 15627.1671
$ws.Range("M56").Value = -5466
$ws.Range("N56").Value = -16695.1671
$ws.Range("H80").Value = 154485.92
$ws.Range("J80").Value = 704.25
$ws.Range("L80").Value = 2112.75
$ws.Range("N80").Value = -4108.75
$ws.Range("H83").Value = 154485.92
$ws.Range("J83").Value = 704.25
$ws.Range("L83").Value = 6338.25
$ws.Range("N83").Value = -16322.25
$ws.Range("H112").Value = 30892.53
$ws.Range("J112").Value = 45045.652
$ws.Range("L112").Value = 135136.956
$ws.Range("N112").Value = -137352.956
$ws.Range("H132").Value = 4388.25
$ws.Range("I132").Value = 3686.875
$ws.Range("J132").Value = 9999.25
$ws.Range("K132").Value = 11060.625
$ws.Range("L132").Value = 29997.75
$ws.Range("M132").Value = -8530.625
$ws.Range("N132").Value = -35057.75
$ws.Range("H137").Value = 6622.0464
$ws.Range("I137").Value = 7225.079
$ws.Range("J137").Value = 2039
$ws.Range("K137").Value = 21675.237
$ws.Range("L137").Value = 6117
$ws.Range("M137").Value = -19125.237
$ws.Range("N137").Value = -11217
$ws.Range("H138").Value = 3655.186
$ws.Range("J138").Value = 4454.516
$ws.Range("L138").Value = 13363.548
$ws.Range("N138").Value = -23643.548

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 42168.85
$ws.Range("I2").Value = 5741.3184
$ws.Range("J2").Value = 202450
$ws.Range("K2").Value = 5741.3184
$ws.Range("L2").Value = 202450
$ws.Range("M2").Value = -5628.3184
$ws.Range("N2").Value = -202676
$ws.Range("H61").Value = 6517.7915
$ws.Range("I61").Value = 7029.231
$ws.Range("J61").Value = 5913.364
$ws.Range("K61").Value = 7029.231
$ws.Range("L61").Value = 5913.364
$ws.Range("M61").Value = -6817.231
$ws.Range("N61").Value = -6337.364
$ws.Range("H74").Value = 4704.923
$ws.Range("I74").Value = 2754.647
$ws.Range("K74").Value = 2754.647
$ws.Range("M74").Value = -1880.647
$ws.Range("H77").Value = 4704.923
$ws.Range("I77").Value = 2754.647
$ws.Range("K77").Value = 13773.235
$ws.Range("M77").Value = -9405.235000000001
$ws.Range("H116").Value = 42168.85
$ws.Range("I116").Value = 5741.3184
$ws.Range("J116").Value = 202450
$ws.Range("K116").Value = 5741.3184
$ws.Range("L116").Value = 202450
$ws.Range("M116").Value = -3447.3184
$ws.Range("N116").Value = -207038
$ws.Range("H122").Value = 393334.56
$ws.Range("I122").Value = 2619.6365
$ws.Range("J122").Value = 1007315.1
$ws.Range("K122").Value = 7858.9095
$ws.Range("L122").Value = 3021945.3
$ws.Range("M122").Value = -5408.9095
$ws.Range("N122").Value = -3026845.3
$ws.Range("H132").Value = 2752.2917
$ws.Range("I132").Value = 1677.1177
$ws.Range("J132").Value = 5363.4287
$ws.Range("K132").Value = 5031.3531
$ws.Range("L132").Value = 16090.2861
$ws.Range("M132").Value = -2501.3531
$ws.Range("N132").Value = -21150.2861
$ws.Range("H136").Value = 6517.7915
$ws.Range("I136").Value = 7029.231
$ws.Range("J136").Value = 5913.364
$ws.Range("K136").Value = 21087.693
$ws.Range("L136").Value = 17740.092
$ws.Range("M136").Value = -18537.693
$ws.Range("N136").Value = -22840.092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 42168.85
$ws.Range("I3").Value = 5741.3184
$ws.Range("J3").Value = 202450
$ws.Range("K3").Value = 5741.3184
$ws.Range("L3").Value = 202450
$ws.Range("M3").Value = -5627.3184
$ws.Range("N3").Value = -202678
$ws.Range("H134").Value = 2320.238
$ws.Range("I134").Value = 1911.5834
$ws.Range("J134").Value = 4772.1665
$ws.Range("K134").Value = 5734.7502
$ws.Range("L134").Value = 14316.4995
$ws.Range("M134").Value = -3199.7502
$ws.Range("N134").Value = -19386.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2134.1667
$ws.Range("I31").Value = 1400.9375
$ws.Range("J31").Value = 8000
$ws.Range("K31").Value = 1400.9375
$ws.Range("L31").Value = 8000
$ws.Range("M31").Value = -1105.9375
$ws.Range("N31").Value = -8590
$ws.Range("H34").Value = 2134.1667
$ws.Range("I34").Value = 1400.9375
$ws.Range("J34").Value = 8000
$ws.Range("K34").Value = 1400.9375
$ws.Range("L34").Value = 8000
$ws.Range("M34").Value = -1198.9375
$ws.Range("N34").Value = -8404
$ws.Range("H74").Value = 72319.664
$ws.Range("J74").Value = 72319.664
$ws.Range("L74").Value = 72319.664
$ws.Range("N74").Value = -74067.664
$ws.Range("H77").Value = 72319.664
$ws.Range("J77").Value = 72319.664
$ws.Range("L77").Value = 216958.992
$ws.Range("N77").Value = -225694.992
$ws.Range("H97").Value = 69999.5
$ws.Range("J97").Value = 79999
$ws.Range("L97").Value = 79999
$ws.Range("N97").Value = -81981
$ws.Range("H134").Value = 4649.5454
$ws.Range("I134").Value = 2526.5334
$ws.Range("J134").Value = 9198.857
$ws.Range("K134").Value = 7579.600199999999
$ws.Range("L134").Value = 27596.571
$ws.Range("M134").Value = -5044.600199999999
$ws.Range("N134").Value = -32666.571
$ws.Range("H139").Value = 55496.25
$ws.Range("J139").Value = 55496.25
$ws.Range("L139").Value = 55496.25
$ws.Range("N139").Value = -65776.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 13889130
$ws.Range("I23").Value = 60.5
$ws.Range("K23").Value = 181.5
$ws.Range("M23").Value = 53.5
$ws.Range("H47").Value = 932.3333
$ws.Range("J47").Value = 932.3333
$ws.Range("L47").Value = 2796.9999
$ws.Range("N47").Value = -3658.9999
$ws.Range("H50").Value = 838.8570999999999
$ws.Range("I50").Value = 645.3333
$ws.Range("J50").Value = 2000
$ws.Range("K50").Value = 1935.9999
$ws.Range("L50").Value = 6000
$ws.Range("M50").Value = -1454.9999
$ws.Range("N50").Value = -6962
$ws.Range("H53").Value = 838.8570999999999
$ws.Range("I53").Value = 645.3333
$ws.Range("J53").Value = 2000
$ws.Range("K53").Value = 1935.9999
$ws.Range("L53").Value = 6000
$ws.Range("M53").Value = -1454.9999
$ws.Range("N53").Value = -6962
$ws.Range("H55").Value = 6782.615
$ws.Range("I55").Value = 1604
$ws.Range("J55").Value = 7724.1816
$ws.Range("K55").Value = 4812
$ws.Range("L55").Value = 23172.5448
$ws.Range("M55").Value = -4635
$ws.Range("N55").Value = -23526.5448
$ws.Range("H131").Value = 1730.0741
$ws.Range("I131").Value = 783.9
$ws.Range("J131").Value = 2286.647
$ws.Range("K131").Value = 2351.7
$ws.Range("L131").Value = 6859.941
$ws.Range("M131").Value = 2688.3
$ws.Range("N131").Value = -16939.941
$ws.Range("H140").Value = 11303.25
$ws.Range("I140").Value = 15235.77
$ws.Range("K140").Value = 45707.31
$ws.Range("M140").Value = -40527.31

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4500
$ws.Range("I113").Value = 4500
$ws.Range("K113").Value = 4500
$ws.Range("M113").Value = -2330
$ws.Range("H122").Value = 11819.375
$ws.Range("J122").Value = 16114.667
$ws.Range("L122").Value = 48344.001
$ws.Range("N122").Value = -53244.001
$ws.Range("H132").Value = 3659
$ws.Range("I132").Value = 3633.8235
$ws.Range("J132").Value = 3730.3333
$ws.Range("K132").Value = 10901.4705
$ws.Range("L132").Value = 11190.9999
$ws.Range("M132").Value = -8371.470499999999
$ws.Range("N132").Value = -16250.9999
$ws.Range("H135").Value = 75691.25
$ws.Range("J135").Value = 75691.25
$ws.Range("L135").Value = 75691.25
$ws.Range("N135").Value = -85831.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7273.28
$ws.Range("I61").Value = 5606
$ws.Range("J61").Value = 19500
$ws.Range("K61").Value = 5606
$ws.Range("L61").Value = 19500
$ws.Range("M61").Value = -5404
$ws.Range("N61").Value = -19904
$ws.Range("H82").Value = 3107.5386
$ws.Range("I82").Value = 4779.8
$ws.Range("J82").Value = 2062.375
$ws.Range("K82").Value = 4779.8
$ws.Range("L82").Value = 2062.375
$ws.Range("M82").Value = -4418.8
$ws.Range("N82").Value = -2784.375
$ws.Range("H85").Value = 3107.5386
$ws.Range("I85").Value = 4779.8
$ws.Range("J85").Value = 2062.375
$ws.Range("K85").Value = 4779.8
$ws.Range("L85").Value = 2062.375
$ws.Range("M85").Value = -3531.8
$ws.Range("N85").Value = -4558.375
$ws.Range("H93").Value = 7334.737
$ws.Range("I93").Value = 9054.929
$ws.Range("J93").Value = 2518.2
$ws.Range("K93").Value = 9054.929
$ws.Range("L93").Value = 2518.2
$ws.Range("M93").Value = -7806.929
$ws.Range("N93").Value = -5014.2
$ws.Range("H100").Value = 114999.5
$ws.Range("I100").Value = 114999.5
$ws.Range("K100").Value = 114999.5
$ws.Range("M100").Value = -114458.5
$ws.Range("H113").Value = 7273.28
$ws.Range("I113").Value = 5606
$ws.Range("J113").Value = 19500
$ws.Range("K113").Value = 5606
$ws.Range("L113").Value = 19500
$ws.Range("M113").Value = -3436
$ws.Range("N113").Value = -23840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5271.28
$ws.Range("I122").Value = 4134.909
$ws.Range("K122").Value = 12404.727
$ws.Range("M122").Value = -9954.726999999999
$ws.Range("H135").Value = 8431134
$ws.Range("J135").Value = 8431134
$ws.Range("L135").Value = 8431134
$ws.Range("N135").Value = -8441274
